$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: FSN, Name, Age, Phone, Vaccine_Dose
$ws.Range("A1").Value = "FSN"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Age"
$ws.Range("D1").Value = "Phone"
$ws.Range("E1").Value = "Vaccine_Dose"

# Widen the Vaccine_Dose column so the header isn't truncated
$ws.Columns.Item(5).ColumnWidth = 14.1796875

# Leave the selection where the author's cursor ended up after entering the headers
[void]$ws.Range("F4").Select()
